# lander/circular_orbit_altitudes.xlsx -- "Deleted commented code, tested all
# scenarios in different situations, autopilot controller now has circular
# checkbox"
#
# The workbook has a single worksheet ("Sheet1") holding the raw altitude
# samples used by a scatter chart + linear trendline ("Chart 2"). The data
# set was extended from 5 rows / 2 columns (A:B) to 18 rows / 3 columns
# (A:B:C, where C is a scenario/run marker: 5, 1, 3), and the chart was
# resized/repositioned and the trendline now also shows R^2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Replace the sample data in A1:C18 (was A1:B5).
# ---------------------------------------------------------------------
$data = @(
    @(5500000, 6537910, 5),
    @(6000000, 7395880, 5),
    @(6500000, 8286530, 5),
    @(7000000, 9182160, 5),
    @(7500000, 10107600, 5),
    @(3900000, 4037580, 5),
    @(5500000, 6556470, 1),
    @(6000000, 7420520, 1),
    @(6500000, 8309000, 1),
    @(7000000, 9215470, 1),
    @(7500000, 10133100, 1),
    @(3900000, 4062140, 1),
    @(5500000, 6703850, 3),
    @(6000000, 7565600, 3),
    @(6500000, 8455870, 3),
    @(7000000, 9362670, 3),
    @(7500000, 10279100, 3),
    @(3900000, 4178350, 3)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

# Selection moved from L13 to L10.
$ws.Range("L10").Select()

# ---------------------------------------------------------------------
# 2. Chart tweaks: the series is no longer smoothed, and the trendline now
#    also displays R^2 (it already displayed the equation).
# ---------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection(1)
$ser.Smooth = $false

$tl = $ser.Trendlines().Item(1)
$tl.DisplayRSquared = $true
$tl.DisplayEquation = $true

# ---------------------------------------------------------------------
# 3. Reposition/resize the chart to its new anchor cells:
#    from D1 (col 3, 203200 EMU / col 0, 139700 EMU)
#    to   P34 (col 15, 546100 EMU / row 33, 50800 EMU)
# ---------------------------------------------------------------------
$co.Left = 191.3125
$co.Top = 11
$co.Width = 728.25
$co.Height = 521
